$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4166666666666667
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = 0.4545454545454545

# Row 3
$ws.Range("B3").Value = 0.5833333333333334
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 0.5384615384615384

# Row 4
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.5

# Row 5
$ws.Range("B5").Value = 0.5
$ws.Range("C5").Value = 0.5
$ws.Range("D5").Value = 0.4965034965034965

# Row 6
$ws.Range("B6").Value = 0.513888888888889
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 0.5034965034965034

# Row 7
$ws.Range("B7").Value = 0.1428571428571428
$ws.Range("C7").Value = 0.1
$ws.Range("D7").Value = 0.1176470588235294

# Row 8
$ws.Range("B8").Value = 0.4705882352941176
$ws.Range("C8").Value = 0.5714285714285714
$ws.Range("D8").Value = 0.5161290322580646

# Row 9
$ws.Range("B9").Value = 0.375
$ws.Range("C9").Value = 0.375
$ws.Range("D9").Value = 0.375
$ws.Range("E9").Value = 0.375

# Row 10
$ws.Range("B10").Value = 0.3067226890756303
$ws.Range("C10").Value = 0.3357142857142857
$ws.Range("D10").Value = 0.316888045540797

# Row 11
$ws.Range("B11").Value = 0.3340336134453781
$ws.Range("C11").Value = 0.375
$ws.Range("D11").Value = 0.3500948766603416

# Row 12
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0

# Row 13
$ws.Range("B13").Value = 0.5833333333333334
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0.7368421052631579

# Row 15
$ws.Range("B15").Value = 0.2916666666666667
$ws.Range("C15").Value = 0.5
$ws.Range("D15").Value = 0.3684210526315789

# Row 16
$ws.Range("B16").Value = 0.3402777777777778
$ws.Range("C16").Value = 0.5833333333333334
$ws.Range("D16").Value = 0.4298245614035088

# Row 22
$ws.Range("B22").Value = 0.4166666666666667
$ws.Range("C22").Value = 0.5
$ws.Range("D22").Value = 0.4545454545454545

# Row 23
$ws.Range("B23").Value = 0.5833333333333334
$ws.Range("C23").Value = 0.5
$ws.Range("D23").Value = 0.5384615384615384

# Row 24
$ws.Range("B24").Value = 0.5
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = 0.5
$ws.Range("E24").Value = 0.5

# Row 25
$ws.Range("B25").Value = 0.5
$ws.Range("C25").Value = 0.5
$ws.Range("D25").Value = 0.4965034965034965

# Row 26
$ws.Range("B26").Value = 0.513888888888889
$ws.Range("C26").Value = 0.5
$ws.Range("D26").Value = 0.5034965034965034
